# Updated cryptos list on Fri Jun  9 10:24:47 UTC 2023 with GitHub Actions
#
# Refreshes the per-coin Price (column D) and Volume(1h) (column E) snapshot
# values on Sheet1, plus the Elrond/Cronos row swap (rows 49-50) that came
# from the source ranking re-sorting between runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. thousand-dot-separated
# prices like "26.643.97", or plain decimals like "1.000") that must stay
# TEXT, matching the workbook source data. A leading apostrophe forces
# Excel to keep them as text instead of auto-coercing to a number.
function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
}

Set-TextValue "D2" "26.643.97"
$ws.Range("E2").Value = "  +0.88%  "

Set-TextValue "D3" "1.844.62"
$ws.Range("E3").Value = "  +0.24%  "

Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue "D5" "259.82"
$ws.Range("E5").Value = "  -0.60%  "

Set-TextValue "D6" "1.000"
$ws.Range("E6").Value = "  -0.01%  "

Set-TextValue "D7" "0.5284"
$ws.Range("E7").Value = "  +1.77%  "

$ws.Range("E8").Value = "  -3.32%  "

Set-TextValue "D9" "0.06805"
$ws.Range("E9").Value = "  +0.43%  "

Set-TextValue "D10" "18.96"
$ws.Range("E10").Value = "  +1.90%  "

Set-TextValue "D11" "0.7866"
$ws.Range("E11").Value = "  +1.41%  "

Set-TextValue "D12" "0.07787"
$ws.Range("E12").Value = "  +0.29%  "

Set-TextValue "D13" "1.851.62"
$ws.Range("E13").Value = "  +0.86%  "

Set-TextValue "D14" "88.26"

Set-TextValue "D15" "5.019"
$ws.Range("E15").Value = "  +0.47%  "

Set-TextValue "D16" "1.001"
$ws.Range("E16").Value = "  +0.12%  "

Set-TextValue "D17" "13.91"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("E18").Value = "  +0.02%  "

Set-TextValue "D19" "0.000007933"
$ws.Range("E19").Value = "  -0.11%  "

Set-TextValue "D20" "26.654.44"
$ws.Range("E20").Value = "  +0.87%  "

Set-TextValue "D21" "2.084.84"
$ws.Range("E21").Value = "  +0.74%  "

Set-TextValue "D22" "4.616"
$ws.Range("E22").Value = "  -0.12%  "

Set-TextValue "D23" "5.986"
$ws.Range("E23").Value = "  +0.18%  "

Set-TextValue "D24" "9.332"
$ws.Range("E24").Value = "  -1.96%  "

Set-TextValue "D25" "143.00"
$ws.Range("E25").Value = "  -1.49%  "

Set-TextValue "D26" "2.225"
$ws.Range("E26").Value = "  +1.81%  "

Set-TextValue "D27" "1.680"
$ws.Range("E27").Value = "  +1.68%  "

Set-TextValue "D28" "17.04"
$ws.Range("E28").Value = "  +0.51%  "

Set-TextValue "D29" "111.00"
$ws.Range("E29").Value = "  -0.56%  "

Set-TextValue "D30" "4.211"
$ws.Range("E30").Value = "  +0.75%  "

Set-TextValue "D31" "0.08719"
$ws.Range("E31").Value = "  +0.21%  "

Set-TextValue "D32" "4.095"
$ws.Range("E32").Value = "  -0.38%  "

Set-TextValue "D33" "0.04897"
$ws.Range("E33").Value = "  +1.73%  "

Set-TextValue "D34" "0.7331"
$ws.Range("E34").Value = "  +2.12%  "

$ws.Range("E35").Value = "  +1.46%  "

Set-TextValue "D36" "2.858"
$ws.Range("E36").Value = "  +0.36%  "

Set-TextValue "D37" "3.109"
$ws.Range("E37").Value = "  +0.77%  "

Set-TextValue "D38" "2.295"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("E39").Value = "  -2.53%  "

Set-TextValue "D40" "0.4830"
$ws.Range("E40").Value = "  +0.31%  "

Set-TextValue "D41" "0.9009"
$ws.Range("E41").Value = "  +0.09%  "

Set-TextValue "D42" "109.92"
$ws.Range("E42").Value = "  -1.29%  "

Set-TextValue "D43" "5.962"
$ws.Range("E43").Value = "  -1.49%  "

Set-TextValue "D44" "1.001"
$ws.Range("E44").Value = "  +0.04%  "

Set-TextValue "D45" "7.720"
$ws.Range("E45").Value = "  +0.11%  "

Set-TextValue "D46" "0.4206"
$ws.Range("E46").Value = "  +1.36%  "

Set-TextValue "D47" "9.091"
$ws.Range("E47").Value = "  +0.98%  "

Set-TextValue "D48" "0.1242"
$ws.Range("E48").Value = "  +1.52%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.05820"
$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D50" "34.85"
$ws.Range("E50").Value = "  -0.41%  "

Set-TextValue "D51" "0.8955"
$ws.Range("E51").Value = "  +1.42%  "

